$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Add holiday marker values (-1) to column B for rows 7 and 13
$ws.Range("B7").Value = -1
$ws.Range("B13").Value = -1

# Update the active selection to reflect the latest editing position
$ws.Range("F11").Select()
